# Vega Central Mapocho de Santiago - Berenjena: insert a new daily record
# as the first data row (row 125), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125 (existing rows 125..233 shift to 126..234)
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(125, 1).Value  = 9
$ws.Cells.Item(125, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(125, 3).Value  = "Metropolitana"
$ws.Cells.Item(125, 4).Value  = 44586
$ws.Cells.Item(125, 5).Value  = 13
$ws.Cells.Item(125, 6).Value  = 100112001
$ws.Cells.Item(125, 7).Value  = "Berenjena"
$ws.Cells.Item(125, 8).Value  = "Sin especificar"
$ws.Cells.Item(125, 9).Value  = "Primera"
$ws.Cells.Item(125, 10).Value = 61
$ws.Cells.Item(125, 11).Value = 10000
$ws.Cells.Item(125, 12).Value = 10000
$ws.Cells.Item(125, 13).Value = 10000
$ws.Cells.Item(125, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(125, 15).Value = "Región Metropolitana"
$ws.Cells.Item(125, 16).Value = 200
$ws.Cells.Item(125, 17).Value = 50
$ws.Cells.Item(125, 18).Value = "Hortaliza"
